# Refresh the Price (D) and Volume(1h) (E) columns with the latest crypto figures.
# D/E are plain-text columns (inlineStr) in the source sheet, so any replacement
# value that reads as a pure number (e.g. '143.80', '1.00') is written through a
# NumberFormat="@" guard and then ClearFormats() so COM keeps it as text (preserving
# trailing zeros) without leaving a permanent style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.711.37'
$ws.Range('E2').Value = '  -1.60%  '

$ws.Range('D3').Value = '2.903.54'
$ws.Range('E3').Value = '  -2.89%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '530.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.80'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.37%  '

$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.06%  '

$ws.Range('D9').Value = '2.910.18'
$ws.Range('E9').Value = '  -2.77%  '

$ws.Range('E10').Value = '  -3.17%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.97'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.75%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.28%  '

$ws.Range('D13').Value = '3.411.02'
$ws.Range('E13').Value = '  -2.72%  '

$ws.Range('D15').Value = '60.635.18'
$ws.Range('E15').Value = '  -1.84%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.79'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.55%  '

$ws.Range('D17').Value = '2.904.93'

$ws.Range('E18').Value = '  -3.05%  '

$ws.Range('E19').Value = '  -1.40%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.76'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.43%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.77'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.95%  '

$ws.Range('E22').Value = '  +0.59%  '

$ws.Range('E23').Value = '  -0.13%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.45'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.87%  '

$ws.Range('D25').Value = '3.019.67'
$ws.Range('E25').Value = '  -2.90%  '

$ws.Range('E26').Value = '  -3.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.180'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.83'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.42%  '

$ws.Range('D30').Value = '0.0₃0869'
$ws.Range('E30').Value = '  -8.05%  '

$ws.Range('E31').Value = '  -0.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.68'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.39%  '

$ws.Range('E33').Value = '  -3.20%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '151.07'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.72%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.37'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.08%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.60'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.99%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.70%  '

$ws.Range('E38').Value = '  -5.61%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.71'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.39%  '

$ws.Range('E40').Value = '  -5.09%  '

$ws.Range('D41').Value = '2.323.98'
$ws.Range('E41').Value = '  -4.67%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.71'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.64%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.646'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0582'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.75'
$ws.Range('D45').ClearFormats()

$ws.Range('E46').Value = '  -0.10%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.96'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.51%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0237'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.94%  '

$ws.Range('E49').Value = '  -2.29%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.35'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.56'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.03%  '
